$d = $word.ActiveDocument

# Change 1: ${ubicacion} -> ${puestoNuevo.gerenciaUbicacion}, and
#           ${incorporacion.nombreDiaDeIncorporacion} -> ${incorporacion.nombreDiaIncorporacion}
$d.Content.Find.Execute(
    "ciudad de `${ubicacion} a Hrs. 08:30 del día `${incorporacion.nombreDiaDeIncorporacion} `${",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ciudad de `${puestoNuevo.gerenciaUbicacion} a Hrs. 08:30 del día `${incorporacion.nombreDiaIncorporacion} `${",
    2
) | Out-Null

# Change 2: incorporacion.fechaDeIncorporacion -> incorporacion.fechaIncorporacion
$d.Content.Find.Execute(
    "incorporacion.fechaDeIncorporacion}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "incorporacion.fechaIncorporacion}",
    2
) | Out-Null

# Change 3: rewrite the sentence around puesto_nuevo.gerencia / ciudadano / designado
$d.Content.Find.Execute(
    "se hizo presente en las oficinas de la `${puesto_nuevo.gerencia} `${ciudadano} `${persona.nombreCompleto} con C.I.  `${persona.ci} `${persona.exp}., `${designado} mediante Resolución Administrativa de Presidencia N° ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "se hizo presente en las oficinas `${puestoNuevo.gerencia} `${persona.ciudadano} `${persona.nombreCompleto} con C.I.  `${persona.ci} `${persona.exp}., `${persona.designado} mediante Resolución Administrativa de Presidencia N.º ",
    2
) | Out-Null

# Change 4: incorporacion.codigoRAP -> incorporacion.codigoRap
$d.Content.Find.Execute(
    "incorporacion.codigoRAP",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "incorporacion.codigoRap",
    2
) | Out-Null

# Change 5: puesto_nuevo.denominacion -> puestoNuevo.denominacion
$d.Content.Find.Execute(
    "puesto_nuevo.denominacion",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "puestoNuevo.denominacion",
    2
) | Out-Null

# Change 6: ${puesto_nuevo.departamento} / ${puesto_nuevo.gerencia} -> camelCase
$d.Content.Find.Execute(
    "} `${puesto_nuevo.departamento} dependiente `${puesto_nuevo.gerencia} ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "} `${puestoNuevo.departamento} dependiente `${puestoNuevo.gerencia} ",
    2
) | Out-Null

# Change 7: ${puesto_nuevo.item} -> ${puestoNuevo.item}
$d.Content.Find.Execute(
    "`${puesto_nuevo.item}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "`${puestoNuevo.item}",
    2
) | Out-Null

# Change 8: table cell ${puesto_nuevo.gerenciaSinConector} (split across 3 runs)
# collapses into a single run ${puestoNuevo.gerente}
$d.Content.Find.Execute(
    "`${puesto_nuevo.gerenciaSinConector}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "`${puestoNuevo.gerente}",
    2
) | Out-Null
